# Edit the numbered-list text inside the "Isosceles Triangle 3" shape on
# slide 1, splitting each of the three list-item paragraphs into multiple
# runs as produced by a user retyping part of the line:
#   "1. this thing"    -> "1. " + "Changed " + "thing"
#   "2 that thing"     -> "2 "  + "other thing"
#   "3. Another thing" -> "3. " + "Final thing"
#
# We locate the text to replace with TextRange.Find(...) (robust to exact
# character offsets) and then assign .Text on that sub-range, which makes
# PowerPoint split the paragraph's run(s) around the edited span instead of
# touching the rest of the paragraph's runs.
#
# Edits are applied back-to-front (last paragraph first) so that earlier
# Find() results/offsets in the still-unmodified part of the text remain
# valid while later parts of the text are being resized.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Paragraph 4: "3. Another thing" -> "3. Final thing"
$f = $tr.Find("Another thing", 0)
$f.Text = "Final thing"

# Paragraph 3: "2 that thing" -> "2 other thing"
$f = $tr.Find("that thing", 0)
$f.Text = "other thing"

# Paragraph 2: "1. this thing" -> "1. Changed thing"
$f = $tr.Find("this ", 0)
$f.Text = "Changed "
